$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell A2 value (SKU text changes) - leave its formatting untouched
$ws.Range("A2").Value = "013742002836M"

# Add new row A3 with new SKU, copying the formatting used by A1
$ws.Range("A1").Copy()
$ws.Range("A3").Select()
$ws.Paste()
$ws.Range("A3").Value = "013742003321M"

# Move active selection to A2
$ws.Range("A2").Select()
